$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column B (the "Gas scale" update):
# old ranges counted down from 35-45 to 0-5 in steps of 10 starting at 35;
# new ranges count up from 0-5 to 36-45, with the first bucket 0-5 and
# subsequent buckets offset by 1 minute (6-15, 16-25, 26-35, 36-45).
$ws.Range("B2").Value = "36 to 45 minutes"
$ws.Range("B3").Value = "26 to 35 minutes"
$ws.Range("B4").Value = "16 to 25 minutes"
$ws.Range("B5").Value = "6 to 15 minutes"
$ws.Range("B6").Value = "0 to 5 minutes "

# Move the active selection to B3
$ws.Range("B3").Select()
